$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timestamp
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 13:42"

# Full country data table (re-sorted by Casos totales desc, with updated figures)
$data = @(
    @("Estados Unidos", 85749, 314, 1868, 82577, 2122, 9, 1304),
    @("China", 81340, 55, 74588, 3460, 1034, 5, 3292),
    @("Italia", 80589, 0, 10361, 62013, 3612, 0, 8215),
    @("España", 64059, 6273, 9357, 49844, 4165, 493, 4858),
    @("Alemania", 47373, 3435, 5673, 41415, 23, 18, 285),
    @("Iran", 32332, 2926, 11133, 18821, 2893, 144, 2378),
    @("Francia", 29155, 0, 4948, 22511, 3375, 0, 1696),
    @("Suiza", 11951, 140, 897, 10857, 203, 5, 197),
    @("Reino Unido", 11658, 0, 135, 10945, 163, 0, 578),
    @("Corea del Sur", 9332, 91, 4528, 4665, 59, 8, 139),
    @("Paises Bajos", 7431, 0, 3, 6994, 761, 0, 434),
    @("Austria", 7317, 408, 225, 7034, 128, 9, 58),
    @("Belgica", 7284, 1049, 858, 6137, 690, 69, 289),
    @("Portugal", 4268, 724, 43, 4149, 71, 16, 76),
    @("Canada", 4043, 0, 228, 3776, 120, 0, 39),
    @("Noruega", 3669, 297, 6, 3647, 70, 2, 16),
    @("Turquia", 3629, 0, 26, 3528, 136, 0, 75),
    @("Australia", 3180, 130, 170, 2997, 23, 0, 13),
    @("Israel", 3035, 342, 79, 2946, 49, 2, 10),
    @("Brasil", 2988, 3, 6, 2905, 296, 0, 77),
    @("Suecia", 2893, 53, 16, 2800, 214, 0, 77),
    @("Malasia", 2161, 130, 259, 1876, 54, 3, 26),
    @("Chequia", 2062, 137, 10, 2043, 34, 0, 9),
    @("Dinamarca", 2010, 133, 1, 1957, 109, 11, 52),
    @("Irlanda", 1819, 0, 5, 1795, 47, 0, 19),
    @("Luxemburgo", 1453, 0, 6, 1438, 3, 0, 9),
    @("Ecuador", 1403, 0, 3, 1366, 58, 0, 34),
    @("Japon", 1387, 0, 359, 981, 57, 0, 47),
    @("Chile", 1306, 0, 22, 1280, 7, 0, 4),
    @("Rumania", 1292, 263, 115, 1153, 32, 1, 24),
    @("Pakistan", 1252, 51, 23, 1220, 7, 0, 9),
    @("Polonia", 1244, 23, 7, 1221, 3, 0, 16),
    @("Tailandia", 1136, 91, 97, 1034, 11, 1, 5),
    @("Indonesia", 1046, 153, 46, 913, 0, 9, 87),
    @("Finlandia", 1038, 80, 10, 1023, 24, 0, 5),
    @("Rusia", 1036, 196, 45, 988, 8, 0, 3),
    @("Arabia Saudita", 1012, 0, 33, 976, 6, 0, 3),
    @("Sudafrica", 927, 0, 12, 913, 7, 2, 2),
    @("Grecia", 892, 0, 42, 823, 53, 0, 27),
    @("Filipinas", 803, 96, 31, 718, 1, 9, 54),
    @("India", 803, 76, 73, 710, 0, 0, 20),
    @("Islandia", 802, 0, 82, 718, 11, 0, 2),
    @("Crucero", 712, 0, 597, 105, 15, 0, 10),
    @("Singapur", 683, 0, 172, 509, 17, 0, 2),
    @("Panama", 674, 0, 2, 663, 20, 0, 9),
    @("Eslovenia", 632, 70, 10, 613, 14, 3, 9),
    @("Argentina", 589, 0, 72, 505, 0, 0, 12),
    @("Mexico", 585, 110, 4, 573, 1, 2, 8),
    @("Peru", 580, 0, 14, 557, 14, 0, 9),
    @("Estonia", 575, 37, 11, 563, 7, 0, 1),
    @("Croacia", 551, 56, 37, 511, 14, 0, 3),
    @("Catar", 549, 0, 43, 506, 6, 0, 0),
    @("Hong Kong", 518, 64, 111, 403, 5, 0, 4),
    @("Egipto", 495, 0, 102, 369, 0, 0, 24),
    @("Colombia", 491, 0, 8, 477, 0, 0, 6),
    @("Republica Dominicana", 488, 0, 3, 475, 0, 0, 10),
    @("Barein", 466, 8, 227, 235, 1, 0, 4),
    @("Irak", 458, 76, 122, 296, 0, 4, 40),
    @("Serbia", 457, 0, 15, 434, 21, 1, 8),
    @("Libano", 391, 23, 23, 361, 3, 1, 7),
    @("Nueva Zelanda", 368, 0, 37, 331, 1, 0, 0),
    @("Argelia", 367, 0, 29, 313, 0, 0, 25),
    @("Lituania", 345, 46, 1, 339, 2, 1, 5),
    @("Emiratos Arabes Unidos", 333, 0, 52, 279, 2, 0, 2),
    @("Armenia", 329, 39, 28, 300, 6, 0, 1),
    @("Hungria", 300, 39, 34, 256, 6, 0, 10),
    @("Letonia", 280, 36, 1, 279, 0, 0, 0),
    @("Bulgaria", 276, 12, 9, 264, 8, 0, 3),
    @("Marruecos", 275, 0, 8, 256, 1, 0, 11),
    @("Eslovaquia", 269, 43, 2, 267, 1, 0, 0),
    @("Principado de Andorra", 267, 43, 1, 263, 6, 0, 3),
    @("Taiwan", 267, 15, 30, 235, 0, 0, 2),
    @("Uruguay", 238, 0, 0, 238, 3, 0, 0),
    @("Costa Rica", 231, 0, 2, 227, 5, 0, 2),
    @("Bosnia y Herzegovina", 230, 39, 5, 222, 1, 0, 3),
    @("Tunez", 227, 30, 2, 219, 10, 1, 6),
    @("Ucrania", 226, 30, 5, 216, 0, 0, 5),
    @("Kuwait", 225, 17, 57, 168, 11, 0, 0),
    @("Republica de Macedonia", 219, 18, 3, 213, 1, 0, 3),
    @("Jordania", 212, 0, 2, 210, 0, 0, 0),
    @("San Marino", 208, 0, 4, 183, 12, 0, 21),
    @("Albania", 186, 12, 31, 147, 3, 2, 8),
    @("Moldavia", 177, 0, 2, 173, 33, 1, 2),
    @("Vietnam", 163, 10, 20, 143, 3, 0, 0),
    @("Burkina Faso", 152, 0, 10, 135, 0, 0, 7),
    @("Republica de Chipre", 146, 0, 4, 139, 3, 0, 3),
    @("Islas Feroe", 144, 4, 54, 90, 2, 0, 0),
    @("Malta", 139, 5, 2, 137, 1, 0, 0),
    @("Ghana", 136, 4, 1, 131, 1, 0, 4),
    @("Reunion", 135, 0, 1, 134, 0, 0, 0),
    @("Oman", 131, 22, 23, 108, 0, 0, 0),
    @("Kazajistan", 125, 12, 2, 122, 0, 0, 1),
    @("Azerbaiyan", 122, 0, 15, 104, 6, 0, 3),
    @("Senegal", 119, 14, 11, 108, 0, 0, 0),
    @("Brunei", 115, 1, 11, 104, 1, 0, 0),
    @("Venezuela", 107, 0, 31, 75, 2, 0, 1),
    @("Sri Lanka", 106, 0, 7, 99, 5, 0, 0),
    @("Camboya", 99, 1, 11, 88, 1, 0, 0),
    @("Costa de Marfil", 96, 0, 3, 93, 0, 0, 0),
    @("Afganistan", 94, 0, 2, 88, 0, 0, 4),
    @("Bielorrusia", 94, 8, 32, 62, 2, 0, 0),
    @("Estado de Palestina", 91, 5, 17, 73, 0, 0, 1),
    @("Camerun", 88, 13, 2, 84, 0, 1, 2),
    @("Uzbekistan", 83, 8, 5, 77, 8, 1, 1),
    @("Martinica", 81, 0, 0, 80, 12, 0, 1),
    @("Mauricio", 81, 0, 0, 79, 1, 0, 2),
    @("Georgia", 81, 2, 13, 68, 1, 0, 0),
    @("Guadalupe", 73, 0, 0, 72, 4, 0, 1),
    @("Montenegro", 70, 1, 0, 69, 1, 0, 1),
    @("Honduras", 68, 1, 0, 67, 0, 0, 1),
    @("Cuba", 67, 0, 1, 64, 2, 0, 2),
    @("Trinidad yTobago", 65, 0, 0, 64, 0, 0, 1),
    @("Nigeria", 65, 0, 3, 61, 0, 0, 1),
    @("Bolivia", 61, 0, 0, 61, 0, 0, 0),
    @("Kirguistan", 58, 14, 0, 58, 0, 0, 0),
    @("Liechtenstein", 56, 0, 0, 56, 0, 0, 0),
    @("Paraguay", 52, 11, 1, 48, 1, 0, 3),
    @("Consejo Danes para los Refugiados", 51, 0, 2, 46, 0, 0, 3),
    @("Ruanda", 50, 0, 0, 50, 0, 0, 0),
    @("Mayotte", 50, 14, 0, 50, 0, 0, 0),
    @("Banglades", 48, 4, 11, 32, 1, 0, 5),
    @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
    @("Gibraltar", 35, 0, 13, 22, 0, 0, 0),
    @("Monaco", 33, 0, 1, 32, 0, 0, 0),
    @("Macao", 33, 0, 10, 23, 0, 0, 0),
    @("Guam", 32, 0, 0, 31, 0, 0, 1),
    @("Kenia", 31, 0, 1, 29, 0, 0, 1),
    @("Polinesia Francesa", 30, 0, 0, 30, 0, 0, 0),
    @("Isla de Man", 29, 3, 0, 29, 0, 0, 0),
    @("Aruba", 28, 0, 1, 27, 0, 0, 0),
    @("Guayana Francesa", 28, 0, 6, 22, 0, 0, 0),
    @("Jamaica", 26, 0, 2, 23, 0, 0, 1),
    @("Togo", 25, 1, 1, 24, 0, 0, 0),
    @("Guatemala", 25, 0, 4, 20, 0, 0, 1),
    @("Madagascar", 24, 1, 0, 24, 0, 0, 0),
    @("Barbados", 24, 0, 0, 24, 0, 0, 0),
    @("Uganda", 18, 4, 0, 18, 0, 0, 0),
    @("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
    @("Zambia", 16, 0, 0, 16, 0, 0, 0),
    @("Etiopia", 16, 4, 0, 16, 0, 0, 0),
    @("Nueva Caledonia", 15, 1, 0, 15, 0, 0, 0),
    @("Bermudas", 15, 0, 2, 13, 0, 0, 0),
    @("Maldivas", 14, 1, 9, 5, 0, 0, 0),
    @("Tanzania", 13, 0, 0, 13, 0, 0, 0),
    @("El Salvador", 13, 0, 0, 13, 0, 0, 0),
    @("Republica de Yibuti", 12, 1, 0, 12, 0, 0, 0),
    @("Guinea Ecuatorial", 12, 0, 0, 12, 0, 0, 0),
    @("Mongolia", 11, 0, 0, 11, 0, 0, 0),
    @("San Martin (Parte Francesa)", 11, 0, 0, 11, 0, 0, 0),
    @("Dominica", 11, 0, 0, 11, 0, 0, 0),
    @("Mali", 11, 7, 0, 11, 0, 0, 0),
    @("Niger", 10, 0, 0, 9, 0, 0, 1),
    @("Bahamas", 9, 0, 1, 8, 0, 0, 0),
    @("Groenlandia", 9, 3, 2, 7, 0, 0, 0),
    @("Haiti", 8, 0, 0, 8, 0, 0, 0),
    @("Guinea", 8, 4, 0, 8, 0, 0, 0),
    @("Surinam", 8, 0, 0, 8, 0, 0, 0),
    @("Islas Caimanes", 8, 0, 0, 7, 0, 0, 1),
    @("Namibia", 8, 0, 2, 6, 0, 0, 0),
    @("Mozambique", 7, 0, 0, 7, 0, 0, 0),
    @("Granada", 7, 0, 0, 7, 0, 0, 0),
    @("Antigua y Barbuda", 7, 0, 0, 7, 0, 0, 0),
    @("Seychelles", 7, 0, 0, 7, 0, 0, 0),
    @("Gabon", 7, 0, 0, 6, 0, 0, 1),
    @("Curazao", 7, 0, 2, 4, 0, 0, 1),
    @("Eritrea", 6, 0, 0, 6, 0, 0, 0),
    @("Laos", 6, 0, 0, 6, 0, 0, 0),
    @("Suazilandia", 6, 0, 0, 6, 0, 0, 0),
    @("Benin", 6, 0, 0, 6, 0, 0, 0),
    @("Birmania", 5, 0, 0, 5, 0, 0, 0),
    @("Montserrat", 5, 0, 0, 5, 0, 0, 0),
    @("Siria", 5, 0, 0, 5, 0, 0, 0),
    @("Fiyi", 5, 0, 0, 5, 0, 0, 0),
    @("Guyana", 5, 0, 0, 4, 0, 0, 1),
    @("Zimbabue", 5, 2, 0, 4, 0, 0, 1),
    @("Cabo Verde", 5, 0, 0, 4, 0, 0, 1),
    @("Santa Sede", 4, 0, 0, 4, 0, 0, 0),
    @("Angola", 4, 0, 0, 4, 0, 0, 0),
    @("Congo", 4, 0, 0, 4, 0, 0, 0),
    @("Republica del Chad", 3, 0, 0, 3, 0, 0, 0),
    @("Liberia", 3, 0, 0, 3, 0, 0, 0),
    @("San Bartolome", 3, 0, 0, 3, 0, 0, 0),
    @("Butan", 3, 1, 0, 3, 0, 0, 0),
    @("Mauritania", 3, 0, 0, 3, 0, 0, 0),
    @("Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0),
    @("San Martin (Parte Holandesa)", 3, 0, 0, 3, 0, 0, 0),
    @("Somalia", 3, 1, 0, 3, 0, 0, 0),
    @("Gambia", 3, 0, 0, 2, 0, 0, 1),
    @("Nepal", 3, 0, 1, 2, 0, 0, 0),
    @("Sudan", 3, 0, 0, 2, 0, 0, 1),
    @("Santa Lucia", 3, 0, 1, 2, 0, 0, 0),
    @("Belice", 2, 0, 0, 2, 0, 0, 0),
    @("Anguila", 2, 0, 0, 2, 0, 0, 0),
    @("Islas Turcas y Caicos", 2, 0, 0, 2, 0, 0, 0),
    @("Guinea-Bisau", 2, 0, 0, 2, 0, 0, 0),
    @("San Cristobal y Nieves", 2, 0, 0, 2, 0, 0, 0),
    @("Islas Virgenes Britanicas", 2, 0, 0, 2, 0, 0, 0),
    @("Nicaragua", 2, 0, 0, 1, 0, 1, 1),
    @("Timor Oriental", 1, 0, 0, 1, 0, 0, 0),
    @("Libia", 1, 0, 0, 1, 0, 0, 0),
    @("San Vicente y las Granadinas", 1, 0, 0, 1, 0, 0, 0),
    @("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    for ($j = 1; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

Write-Host "Updated $($data.Count) rows"
